# Apply the employee absence data updates described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data as [employee_id, employee_name, department, absence_reason, absence_duration, absence_date(serial), salary]
$data = @(
    @(2, 19391, "Srta. Hadassa Cavalcanti", "Vendas", "Doenca", 3, 45082, 2864.47),
    @(3, 81762, "Antônio Duarte", "Marketing", "Consulta medica", 8, 45095, 4606.81),
    @(4, 16982, "Ana Carolina Pinto", "Financeiro", "Outros", 7, 45091, 7226.6),
    @(5, 67726, "Liam Farias", "Juridico", "Problemas pessoais", 4, 45101, 7731.44),
    @(6, 93049, "Dra. Joana Farias", "Recursos Humanos", "Consulta medica", 4, 45092, 3148.81),
    @(7, 62117, "Melina Aparecida", "Recursos Humanos", "Doenca", 6, 45093, 7220.21),
    @(8, 99415, "Stella Cardoso", "Operacoes", "Doenca", 6, 45105, 5020.45),
    @(9, 63279, "Vitor Hugo Pimenta", "Atendimento ao Cliente", "Consulta medica", 3, 45094, 5958.24),
    @(10, 57980, "Amanda Costela", "Marketing", "Doenca", 2, 45085, 4954.74),
    @(11, 31790, "Alice Martins", "Atendimento ao Cliente", "Doenca", 1, 45106, 6285.94)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
}
